{"js": "// New text values, in document order: paragraph 0 is the date line,\n// paragraphs 1..100 are the table-cell math problems (row-major, left to\n// right), matching the order Word exposes via body.paragraphs.items.\nconst newTexts = [\n  \"2023-03-05 Sunday\",\n  \"29+54=\", \"52+26=\", \"7+49=\", \"12+56=\", \"34+56=\",\n  \"77-52=\", \"65+0=\", \"35+45=\", \"16+74=\", \"7+48=\",\n  \"88-46=\", \"29+3=\", \"1+56=\", \"70-31=\", \"9+62=\",\n  \"30+48=\", \"84-28=\", \"98-68=\", \"66+33=\", \"74-28=\",\n  \"96-20=\", \"17+63=\", \"57-1=\", \"98-73=\", \"71-18=\",\n  \"77-68=\", \"21-9=\", \"86+11=\", \"37+29=\", \"89-71=\",\n  \"88+8=\", \"2+70=\", \"93-82=\", \"19+37=\", \"25+28=\",\n  \"4+14=\", \"87-74=\", \"57+28=\", \"47-33=\", \"93-13=\",\n  \"85-69=\", \"72-33=\", \"18+3=\", \"61-46=\", \"79-72=\",\n  \"16+33=\", \"92-30=\", \"4+92=\", \"16+81=\", \"95-5=\",\n  \"3+28=\", \"26+59=\", \"14+83=\", \"72-33=\", \"89-53=\",\n  \"62+2=\", \"3+70=\", \"88-72=\", \"75-31=\", \"25+56=\",\n  \"22+32=\", \"95-49=\", \"40+50=\", \"92-45=\", \"33-20=\",\n  \"14-8=\", \"57-54=\", \"63+1=\", \"53+42=\", \"76-4=\",\n  \"79-74=\", \"52+18=\", \"95-37=\", \"56-54=\", \"43+11=\",\n  \"51+10=\", \"30-29=\", \"97-88=\", \"14+78=\", \"66-52=\",\n  \"32+10=\", \"98-17=\", \"23+21=\", \"0+38=\", \"42-10=\",\n  \"75-42=\", \"64+19=\", \"67-10=\", \"67+27=\", \"88-44=\",\n  \"38+16=\", \"63-17=\", \"93-92=\", \"10+76=\", \"50-47=\",\n  \"31-15=\", \"25+0=\", \"14+67=\", \"9+67=\", \"67+26=\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst count = Math.min(paragraphs.items.length, newTexts.length);\nfor (let i = 0; i < count; i++) {\n  const para = paragraphs.items[i];\n  const newText = newTexts[i];\n  // Only touch paragraphs whose text actually changes, replacing the whole\n  // paragraph's contents (keeps the existing run formatting since\n  // insertText(\"Replace\") overwrites the paragraph's range in place).\n  if (para.text !== newText) {\n    para.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ Old = '2023-03-04 Saturday'; New = '2023-03-05 Sunday' },\n    @{ Old = '20+12='; New = '29+54=' },\n    @{ Old = '46+10='; New = '52+26=' },\n    @{ Old = '62+22='; New = '7+49=' },\n    @{ Old = '82-20='; New = '12+56=' },\n    @{ Old = '98-67='; New = '34+56=' },\n    @{ Old = '46-35='; New = '77-52=' },\n    @{ Old = '10+84='; New = '65+0=' },\n    @{ Old = '65-37='; New = '35+45=' },\n    @{ Old = '97-95='; New = '16+74=' },\n    @{ Old = '2+68='; New = '7+48=' },\n    @{ Old = '50-35='; New = '88-46=' },\n    @{ Old = '45-29='; New = '29+3=' },\n    @{ Old = '16+6='; New = '1+56=' },\n    @{ Old = '12+80='; New = '70-31=' },\n    @{ Old = '5-2='; New = '9+62=' },\n    @{ Old = '36+54='; New = '30+48=' },\n    @{ Old = '94-1='; New = '84-28=' },\n    @{ Old = '46-28='; New = '98-68=' },\n    @{ Old = '34+58='; New = '66+33=' },\n    @{ Old = '6+93='; New = '74-28=' },\n    @{ Old = '54-32='; New = '96-20=' },\n    @{ Old = '77-51='; New = '17+63=' },\n    @{ Old = '26-16='; New = '57-1=' },\n    @{ Old = '82-7='; New = '98-73=' },\n    @{ Old = '42+36='; New = '71-18=' },\n    @{ Old = '25+3='; New = '77-68=' },\n    @{ Old = '17+11='; New = '21-9=' },\n    @{ Old = '49+12='; New = '86+11=' },\n    @{ Old = '8+21='; New = '37+29=' },\n    @{ Old = '20+75='; New = '89-71=' },\n    @{ Old = '75-75='; New = '88+8=' },\n    @{ Old = '96-73='; New = '2+70=' },\n    @{ Old = '35-13='; New = '93-82=' },\n    @{ Old = '81-29='; New = '19+37=' },\n    @{ Old = '46-42='; New = '25+28=' },\n    @{ Old = '53-5='; New = '4+14=' },\n    @{ Old = '93-16='; New = '87-74=' },\n    @{ Old = '49-39='; New = '57+28=' },\n    @{ Old = '2+4='; New = '47-33=' },\n    @{ Old = '79+1='; New = '93-13=' },\n    @{ Old = '94-76='; New = '85-69=' },\n    @{ Old = '7+50='; New = '72-33=' },\n    @{ Old = '5+15='; New = '18+3=' },\n    @{ Old = '70-20='; New = '61-46=' },\n    @{ Old = '83-58='; New = '79-72=' },\n    @{ Old = '25+74='; New = '16+33=' },\n    @{ Old = '69-0='; New = '92-30=' },\n    @{ Old = '86-66='; New = '4+92=' },\n    @{ Old = '3+33='; New = '16+81=' },\n    @{ Old = '36-28='; New = '95-5=' },\n    @{ Old = '50+7='; New = '3+28=' },\n    @{ Old = '9+18='; New = '26+59=' },\n    @{ Old = '99-10='; New = '14+83=' },\n    @{ Old = '44+2='; New = '72-33=' },\n    @{ Old = '35+31='; New = '89-53=' },\n    @{ Old = '57+41='; New = '62+2=' },\n    @{ Old = '31+18='; New = '3+70=' },\n    @{ Old = '18+32='; New = '88-72=' },\n    @{ Old = '44+1='; New = '75-31=' },\n    @{ Old = '97-61='; New = '25+56=' },\n    @{ Old = '92-73='; New = '22+32=' },\n    @{ Old = '46-13='; New = '95-49=' },\n    @{ Old = '0+63='; New = '40+50=' },\n    @{ Old = '18+47='; New = '92-45=' },\n    @{ Old = '18+72='; New = '33-20=' },\n    @{ Old = '1+77='; New = '14-8=' },\n    @{ Old = '87+12='; New = '57-54=' },\n    @{ Old = '46+27='; New = '63+1=' },\n    @{ Old = '91-84='; New = '53+42=' },\n    @{ Old = '99-46='; New = '76-4=' },\n    @{ Old = '1+84='; New = '79-74=' },\n    @{ Old = '31+19='; New = '52+18=' },\n    @{ Old = '68+3='; New = '95-37=' },\n    @{ Old = '5+11='; New = '56-54=' },\n    @{ Old = '18+73='; New = '43+11=' },\n    @{ Old = '38-16='; New = '51+10=' },\n    @{ Old = '96-54='; New = '30-29=' },\n    @{ Old = '48-8='; New = '97-88=' },\n    @{ Old = '59-39='; New = '14+78=' },\n    @{ Old = '89-74='; New = '66-52=' },\n    @{ Old = '41-37='; New = '32+10=' },\n    @{ Old = '60+24='; New = '98-17=' },\n    @{ Old = '72-40='; New = '23+21=' },\n    @{ Old = '84-7='; New = '0+38=' },\n    @{ Old = '30+11='; New = '42-10=' },\n    @{ Old = '78-39='; New = '75-42=' },\n    @{ Old = '27-23='; New = '64+19=' },\n    @{ Old = '43-22='; New = '67-10=' },\n    @{ Old = '50+12='; New = '67+27=' },\n    @{ Old = '81-57='; New = '88-44=' },\n    @{ Old = '24-11='; New = '38+16=' },\n    @{ Old = '52-24='; New = '63-17=' },\n    @{ Old = '65-36='; New = '93-92=' },\n    @{ Old = '49-18='; New = '10+76=' },\n    @{ Old = '11+68='; New = '50-47=' },\n    @{ Old = '6+7='; New = '31-15=' },\n    @{ Old = '16+51='; New = '25+0=' },\n    @{ Old = '89-61='; New = '14+67=' },\n    @{ Old = '91-21='; New = '9+67=' },\n    @{ Old = '50-20='; New = '67+26=' },\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n\nWrite-Output \"done\""}
